$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update G2 value
$ws.Range("G2").Value = 11

# Add new row 9
$ws.Range("A9").Value = "Aceite de Canola"
$ws.Range("B9").Value = "Aceite"
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0.01
$ws.Range("E9").Value = 884
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 1500

$ws.Range("D9").Select()
